# Commit: "Update all ticket price mentions to 175 AED"
#
# The underlying xlsx diff for this commit is a new booking row appended
# to the bookings sheet (row 37), recorded at the already-standard
# ticket price of 175 AED (which every other row in the sheet already
# uses) - i.e. this particular booking is an instance of that 175 AED
# price being used/confirmed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append after the last used row (so this keeps working even if the
# sheet already grew), but the data set here matches the known new row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row = $lastRow + 1

$ws.Cells.Item($row, 1).Value  = "SL-20251206-001"   # booking_id
$ws.Cells.Item($row, 2).Value  = "2025-12-06 20:06:27" # created_at
$ws.Cells.Item($row, 3).Value  = "fahaf"              # name

# phone column is stored as text for every existing row (even values
# that are all digits, e.g. "0502992932") - force text here too so the
# type matches the rest of the column instead of becoming a number.
$ws.Cells.Item($row, 4).Value  = "'1234"              # phone

$ws.Cells.Item($row, 5).Value  = 1                    # tickets
$ws.Cells.Item($row, 6).Value  = 175                  # ticket_price
$ws.Cells.Item($row, 7).Value  = 175                  # total_amount
$ws.Cells.Item($row, 8).Value  = "pending"             # status
$ws.Cells.Item($row, 10).Value = "pending"             # payment_status
# payment_intent_id, redirect_url, notes are left blank for this booking.
